$d = $word.ActiveDocument
$d.Bookmarks.ShowHidden = $true

# --------------------------------------------------------------------
# 1. Append the new sentence to the paragraph ending "...so their code
#    'escaped' review)." (right after the existing trailing period).
# --------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("review).")
$anchor.Collapse(0)

$newSentence = " One particular example of this from the feedback video was about scrum meeting notes (readme file) were updated directly on GitHub which pushed to the master branch which is not a good practice so those changes will be from now on pushed to a separate branch first then merged in to the master/dev branch."

# Insert the sentence plus a one-character placeholder ("Z") so that the
# insertion point for the relocated _GoBack bookmark is never exactly on
# a paragraph-end boundary (that position is unreliable for Bookmarks.Add).
$anchor.InsertAfter($newSentence + "Z")

# --------------------------------------------------------------------
# 2. Relocate the _GoBack bookmark to the end of the paragraph we just
#    extended (this mirrors Word leaving _GoBack at the last edit spot).
# --------------------------------------------------------------------
$bmPos = $anchor.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Item("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the placeholder character now that the bookmark is anchored.
$placeholder = $d.Range($bmPos, $bmPos + 1)
$placeholder.Delete()

# --------------------------------------------------------------------
# 3. Fix up "The sprint bu" + "rndown chart..." so the word "burndown"
#    reads as a whole word (it used to be split by the old bookmark
#    location) and mark it off as its own run (closest achievable
#    approximation of the spell-checker's proofErr wrapper).
# --------------------------------------------------------------------
$buRange = $d.Content
$buRange.Find.Execute("The sprint bu")
$buEnd = $buRange.End
$bu = $d.Range($buEnd - 2, $buEnd)
$bu.Delete()

$burndownIns = $d.Range($buEnd - 2, $buEnd - 2)
$burndownIns.InsertAfter("burndown")
$burndownIns.Bold = 1
$burndownIns.Bold = 0

$rndownRange = $d.Content
$rndownRange.Find.Execute("rndown chart")
$rndownStart = $rndownRange.Start
$rndown = $d.Range($rndownStart, $rndownStart + 6)
$rndown.Delete()

# The empty bookmark pair left behind by the original "_GoBack" location
# is no longer needed once the word is unified.
$d.Bookmarks.Item("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $bmRange)
